$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two new blank rows right after the existing data rows
#    (row 16 = CC/ERIKA, row 17 = PPT/MINERVA) so they become rows
#    18 and 19.  Everything below (signature block) shifts down by 2.
# ------------------------------------------------------------------
$ws.Rows("18:19").Insert()

# ------------------------------------------------------------------
# 2. Duplicate the two worker rows (16:17) into the freshly inserted
#    rows (18:19), carrying over values, number formats and borders.
# ------------------------------------------------------------------
$ws.Range("B16:J17").Copy($ws.Range("B18:J19"))

# ------------------------------------------------------------------
# 3. The new rows represent a new "Periodo Mora" (2509) instead of
#    the existing 2508, so update column E on the new rows.
# ------------------------------------------------------------------
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2509"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2509"

# ------------------------------------------------------------------
# 4. Center the "Periodo Mora" values for all four worker rows.
# ------------------------------------------------------------------
$ws.Range("E16:E19").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------------
# 5. Update the totals now that there are 2 periods instead of 1:
#    - "Cant. Periodos" (F13) goes from 1 to 2
#    - "VALOR MORA" (E11) doubles from 113880 to 227760
# ------------------------------------------------------------------
$ws.Range("F13").Value = 2
$ws.Range("E11").Value = 227760

